# "tra cuu diem ttth"
#
# The sheet's whole used range (A1:L3) was selected and formatted as Text
# (numFmtId 49, format code "@") instead of the old General/Date look, and
# the "mssv" column (B2/B3) was re-keyed as text so the leading zero in the
# student ids is preserved (123456789 -> 0306161291 / 0306161292). The
# active selection was also left on G11 instead of K5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select/format the whole used range as Text.
$ws.Cells.NumberFormat = "@"

# Re-enter the mssv values as text (format is already Text, so the leading
# zero is kept instead of being parsed back into a number).
$ws.Range("B2").Value = "0306161291"
$ws.Range("B3").Value = "0306161292"

# Leave the selection on G11, like in the saved workbook.
[void]$ws.Range("G11").Select()
